# MENT-213: Create Questions Category on the webApp
# Add two new health-facility rows to the "HF" sheet (Pebane/PS Cutal and
# Namacurra/PS Naciaia, both under province ZAMBEZIA).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HF")

$lastRow = 94

$ws.Range("B" + ($lastRow + 1)).Value = "ZAMBEZIA"
$ws.Range("C" + ($lastRow + 1)).Value = "Pebane"
$ws.Range("D" + ($lastRow + 1)).Value = "PS Cutal"

$ws.Range("B" + ($lastRow + 2)).Value = "ZAMBEZIA"
$ws.Range("C" + ($lastRow + 2)).Value = "Namacurra"
$ws.Range("D" + ($lastRow + 2)).Value = "PS Naciaia"

# Match the author's resulting selection/scroll state after entering the data.
$ws.Range("D84").Select()
